$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Right-align the "N/A"-style date header cell A20 (it already carries the
#     text "16/07/2018" number format; this adds horizontal=right alignment
#     to that same cell's format, same as Excel would do from the UI). ---
$ws.Range("A20").HorizontalAlignment = -4152

# --- Append a new timesheet row (row 21), matching the pattern used by the
#     other "N/A" rows (e.g. row 15/20): real date, two "N/A" text entries,
#     and a start time, with no end time recorded yet. ---
$ws.Range("A15:D15").Copy()
$ws.Range("A21:D21").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A21").Value = 43298
$ws.Range("B21").Value = "N/A"
$ws.Range("C21").Value = "N/A"
$ws.Range("D21").Value = 0.3125

$ws.Range("E21").Select()
